$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells stay text (avoid numeric auto-coercion of dotted values)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.528.17'
$ws.Range("E2").Value = '  +7.14%  '
$ws.Range("D3").Value = '1.726.70'
$ws.Range("E3").Value = '  +3.82%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").Value = '332.05'
$ws.Range("E5").Value = '  +0.41%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("E7").Value = '  +2.29%  '
$ws.Range("D10").Value = '1.182'
$ws.Range("E10").Value = '  +3.69%  '
$ws.Range("D11").Value = '0.07438'
$ws.Range("E11").Value = '  +5.43%  '
$ws.Range("D12").Value = '1.000'
$ws.Range("E12").Value = '  -0.22%  '
$ws.Range("D13").Value = '6.402'
$ws.Range("E13").Value = '  +5.52%  '
$ws.Range("D14").Value = '20.06'
$ws.Range("E14").Value = '  +2.48%  '
$ws.Range("D15").Value = '7.037'
$ws.Range("E15").Value = '  +6.59%  '
$ws.Range("D16").Value = '1.723.41'
$ws.Range("E16").Value = '  +3.62%  '
$ws.Range("D17").Value = '0.00001073'
$ws.Range("E17").Value = '  +2.25%  '
$ws.Range("D18").Value = '0.06657'
$ws.Range("E18").Value = '  +0.70%  '
$ws.Range("D19").Value = '82.04'
$ws.Range("E19").Value = '  +4.34%  '
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.06%  '
$ws.Range("D21").Value = '16.53'
$ws.Range("E21").Value = '  +4.56%  '
$ws.Range("D22").Value = '6.183'
$ws.Range("E22").Value = '  +4.30%  '
$ws.Range("D23").Value = '12.72'
$ws.Range("E23").Value = '  +1.99%  '
$ws.Range("D24").Value = '26.493.19'
$ws.Range("E24").Value = '  +6.88%  '
$ws.Range("D25").Value = '2.452'
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("D26").Value = '1.425'
$ws.Range("E26").Value = '  +21.66%  '
$ws.Range("D27").Value = '2.388'
$ws.Range("E27").Value = '  -0.87%  '
$ws.Range("D28").Value = '150.47'
$ws.Range("E28").Value = '  +1.26%  '
$ws.Range("D29").Value = '19.46'
$ws.Range("E29").Value = '  +4.27%  '
$ws.Range("D30").Value = '1.915.45'
$ws.Range("E30").Value = '  +3.71%  '
$ws.Range("D31").Value = '131.40'
$ws.Range("E31").Value = '  +4.54%  '
$ws.Range("D32").Value = '4.101'
$ws.Range("E32").Value = '  +0.65%  '
$ws.Range("D33").Value = '5.977'
$ws.Range("E33").Value = '  +4.94%  '
$ws.Range("D34").Value = '0.08574'
$ws.Range("E34").Value = '  +1.09%  '
$ws.Range("D35").Value = '1.686'
$ws.Range("E35").Value = '  +3.10%  '
$ws.Range("D36").Value = '12.72'
$ws.Range("E36").Value = '  +4.92%  '
$ws.Range("D37").Value = '5.354'
$ws.Range("E37").Value = '  +3.86%  '
$ws.Range("D38").Value = '0.02337'
$ws.Range("E38").Value = '  +3.39%  '
$ws.Range("D41").Value = '8.397'
$ws.Range("E41").Value = '  +2.31%  '
$ws.Range("D42").Value = '1.217'
$ws.Range("E42").Value = '  -0.63%  '
$ws.Range("D43").Value = '0.6208'
$ws.Range("E43").Value = '  +4.93%  '
$ws.Range("D44").Value = '14.25'
$ws.Range("E44").Value = '  +5.47%  '
$ws.Range("D45").Value = '1.000'
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("D46").Value = '3.898'
$ws.Range("E46").Value = '  +1.21%  '
$ws.Range("D47").Value = '0.6012'
$ws.Range("E47").Value = '  +6.28%  '
$ws.Range("D48").Value = '128.36'
$ws.Range("E48").Value = '  +2.26%  '
$ws.Range("D49").Value = '2.042'
$ws.Range("E49").Value = '  +4.85%  '
$ws.Range("D50").Value = '0.07175'
$ws.Range("E50").Value = '  +2.95%  '
$ws.Range("D51").Value = '76.95'
$ws.Range("E51").Value = '  +2.78%  '

# Rows where coin rank order swapped
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").Value = '0.3393'
$ws.Range("E8").Value = '  +4.67%  '
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").Value = '48.09'
$ws.Range("E9").Value = '  +1.64%  '
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").Value = '0.2157'
$ws.Range("E39").Value = '  +3.68%  '
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = '0.06206'
$ws.Range("E40").Value = '  +2.55%  '

# Restore default style on Price column (NumberFormat change above bumped style index)
$ws.Range("D2:D51").Style = "Normal"
